$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Uday"
$ws.Range("B2").Value = "Sharma"
$ws.Range("C2").Value = "sharmauday1999@gmail.com"
$ws.Range("D2").Value = "uday"
$ws.Range("E2").Value = "uday"
